$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 822
$ws1.Range("F4").Value = 1405
$ws1.Range("F5").Value = 839
$ws1.Range("F7").Value = 628
$ws1.Range("F10").Value = 49
$ws1.Range("F13").Value = 1605
$ws1.Range("F17").Value = 72
$ws1.Range("F23").Value = 730
$ws1.Range("F25").Value = 1454
$ws1.Range("F26").Value = 175

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 652
$ws2.Range("F7").Value = 271

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 822
$ws4.Range("F5").Value = 1405
$ws4.Range("F6").Value = 839
$ws4.Range("F10").Value = 628
$ws4.Range("F11").Value = 652
$ws4.Range("F14").Value = 49
$ws4.Range("F17").Value = 1605
$ws4.Range("F22").Value = 72
$ws4.Range("F26").Value = 271
$ws4.Range("F35").Value = 730
$ws4.Range("F37").Value = 1454
$ws4.Range("F38").Value = 175
